$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 93
$ws.Range("I6").Value = 79.333336
$ws.Range("J6").Value = 154.5
$ws.Range("K6").Value = 238.000008
$ws.Range("L6").Value = 463.5
$ws.Range("M6").Value = -126.000008
$ws.Range("N6").Value = -687.5

$ws.Range("H28").Value = 2955.6
$ws.Range("I28").Value = 1470.4166
$ws.Range("K28").Value = 1470.4166
$ws.Range("M28").Value = -985.4166

$ws.Range("H32").Value = 539.8
$ws.Range("J32").Value = 566.3333
$ws.Range("L32").Value = 566.3333
$ws.Range("N32").Value = -1218.3333

$ws.Range("H57").Value = 38999
$ws.Range("I57").Value = 38999
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 116997
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -116498
$ws.Range("N57").ClearContents()

$ws.Range("H69").Value = 7074.5
$ws.Range("J69").Value = 7074.5
$ws.Range("L69").Value = 21223.5
$ws.Range("N69").Value = -22971.5

$ws.Range("H72").Value = 7074.5
$ws.Range("J72").Value = 7074.5
$ws.Range("L72").Value = 63670.5
$ws.Range("N72").Value = -72406.5

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws.Range("H111").Value = 1807.7727
$ws.Range("I111").Value = 876
$ws.Range("J111").Value = 1954.8948
$ws.Range("K111").Value = 2628
$ws.Range("L111").Value = 5864.6844
$ws.Range("M111").Value = 439
$ws.Range("N111").Value = -11998.6844

$ws.Range("H138").Value = 4494.5386
$ws.Range("J138").Value = 4729.75
$ws.Range("L138").Value = 14189.25
$ws.Range("N138").Value = -24469.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4704
$ws.Range("I61").Value = 3058.9
$ws.Range("K61").Value = 3058.9
$ws.Range("M61").Value = -2846.9

$ws.Range("H74").Value = 1873.7693
$ws.Range("I74").Value = 1873.7693
$ws.Range("K74").Value = 1873.7693
$ws.Range("M74").Value = -999.7692999999999

$ws.Range("H77").Value = 1873.7693
$ws.Range("I77").Value = 1873.7693
$ws.Range("K77").Value = 9368.8465
$ws.Range("M77").Value = -5000.8465

$ws.Range("H136").Value = 4704
$ws.Range("I136").Value = 3058.9
$ws.Range("K136").Value = 9176.700000000001
$ws.Range("M136").Value = -6626.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 145000
$ws.Range("J74").Value = 145000
$ws.Range("L74").Value = 145000
$ws.Range("N74").Value = -146872

$ws.Range("H77").Value = 145000
$ws.Range("J77").Value = 145000
$ws.Range("L77").Value = 435000
$ws.Range("N77").Value = -444360

$ws.Range("H92").Value = 166666.33
$ws.Range("J92").Value = 166666.33
$ws.Range("L92").Value = 166666.33
$ws.Range("N92").Value = -171658.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1567.1428
$ws.Range("J58").Value = 2613.2856
$ws.Range("L58").Value = 2613.2856
$ws.Range("N58").Value = -3019.2856

$ws.Range("H99").Value = 2199
$ws.Range("I99").Value = 2123.75
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2123.75
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -625.75
$ws.Range("N99").Value = -5496

$ws.Range("H105").Value = 1378.35
$ws.Range("I105").Value = 823.3333
$ws.Range("K105").Value = 823.3333
$ws.Range("M105").Value = 923.6667

$ws.Range("H126").Value = 2199
$ws.Range("I126").Value = 2123.75
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 6371.25
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -3901.25
$ws.Range("N126").Value = -12440

$ws.Range("H136").Value = 1567.1428
$ws.Range("J136").Value = 2613.2856
$ws.Range("L136").Value = 7839.8568
$ws.Range("N136").Value = -12939.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 7560.6665
$ws.Range("I136").Value = 6989
$ws.Range("K136").Value = 20967
$ws.Range("M136").Value = -15867

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H122").Value = 1002897.4
$ws.Range("I122").Value = 5000000
$ws.Range("K122").Value = 15000000
$ws.Range("M122").Value = -14997550

$ws.Range("H136").Value = 20000
$ws.Range("J136").Value = 20000
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -65100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7714
$ws.Range("I7").Value = 7714
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 7714
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -7602
$ws.Range("N7").ClearContents()

$ws.Range("H55").Value = 1078.0667
$ws.Range("J55").Value = 1399.8
$ws.Range("L55").Value = 1399.8
$ws.Range("N55").Value = -1745.8

$ws.Range("H87").Value = 56797
$ws.Range("J87").Value = 56797
$ws.Range("L87").Value = 56797
$ws.Range("N87").Value = -59043

$ws.Range("H90").Value = 56797
$ws.Range("J90").Value = 56797
$ws.Range("L90").Value = 170391
$ws.Range("N90").Value = -181623

$ws.Range("H126").Value = 7714
$ws.Range("I126").Value = 7714
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 23142
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -20672
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 13748.75
$ws.Range("I132").Value = 9667
$ws.Range("K132").Value = 29001
$ws.Range("M132").Value = -26471

$ws.Range("H136").Value = 3107.739
$ws.Range("I136").Value = 2377.3076
$ws.Range("J136").Value = 4057.3
$ws.Range("K136").Value = 7131.9228
$ws.Range("L136").Value = 12171.9
$ws.Range("M136").Value = -4581.9228
$ws.Range("N136").Value = -17271.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2023
$ws.Range("I107").Value = 2265.8333
$ws.Range("K107").Value = 6797.499899999999
$ws.Range("M107").Value = -4877.499899999999

$ws.Range("H126").Value = 5182.6665
$ws.Range("I126").Value = 3315.5557
$ws.Range("J126").Value = 7983.3335
$ws.Range("K126").Value = 9946.667099999999
$ws.Range("L126").Value = 23950.0005
$ws.Range("M126").Value = -7476.667099999999
$ws.Range("N126").Value = -28890.0005

$ws.Range("H132").Value = 1809.9375
$ws.Range("I132").Value = 1618.5
$ws.Range("J132").Value = 3150
$ws.Range("K132").Value = 4855.5
$ws.Range("L132").Value = 9450
$ws.Range("M132").Value = -2325.5
$ws.Range("N132").Value = -14510

$ws.Range("H136").Value = 3583.8147
$ws.Range("I136").Value = 2748.077
$ws.Range("K136").Value = 8244.231
$ws.Range("M136").Value = -5694.231
